# Lattice-multiplication worksheet refresh: every exercise cell in the
# single 5x3 table gets a new "A x B" problem (same visual template -
# sz 32 run, a title line, the second factor split across two digits,
# a "----" rule, and the first factor's two digits used as the left-hand
# lattice labels).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# New (A, B) pairs, in row-major order (row1: col1..col3, row2: col1..col3, ...)
$pairs = @(
    @(78,63), @(25,39), @(52,58),
    @(30,67), @(32,24), @(55,43),
    @(48,17), @(21,14), @(53,40),
    @(17,24), @(80,20), @(96,55),
    @(44,82), @(16,36), @(12,48)
)

$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $pair = $pairs[$idx]
        $idx++
        $aStr = [string]$pair[0]
        $bStr = [string]$pair[1]
        $a0 = $aStr.Substring(0,1)
        $a1 = $aStr.Substring(1,1)
        $b0 = $bStr.Substring(0,1)
        $b1 = $bStr.Substring(1,1)

        # Use -f (format) throughout: this interpreter's "+" on two
        # digit-only strings does numeric addition, not concatenation.
        $line1 = "{0} x {1}" -f $aStr, $bStr
        $line2 = "  {0}    {1}" -f $b0, $b1
        $line3 = "  ----"
        $line4 = "{0}|    |" -f $a0
        $line5 = "{0}|    |" -f $a1

        $xml = "<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
               "<w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr>" +
               ("<w:t>{0}</w:t><w:br/>" -f $line1) +
               ("<w:t xml:space=`"preserve`">{0}</w:t><w:br/>" -f $line2) +
               ("<w:t xml:space=`"preserve`">{0}</w:t><w:br/>" -f $line3) +
               ("<w:t>{0}</w:t><w:br/>" -f $line4) +
               ("<w:t>{0}</w:t>" -f $line5) +
               "</w:r></w:p>"

        $cell = $t.Cell($r, $c)
        $cell.Range.InsertXML($xml)
    }
}
